$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Class" column (H), and the "Password"/"Re-Password" columns (N, O).
# Delete from right to left so earlier deletions don't shift the column
# letters of the ones still queued for deletion.
$ws.Range("N1:O1").EntireColumn.Delete()
$ws.Range("H1:H1").EntireColumn.Delete()

# The Gender list validation (was on C, driven by $Q$1:$R$1) and the
# Designation list validation (was on I, driven by $AE$1:$AF$1) keep their
# old formula references after the column delete, so repoint them at the
# cells that now hold "Male"/"Female" and "Faculty"/"Asst. Faculty".
$ws.Range("C1:C1048576").Validation.Modify(3, 1, 1, "=`$N`$1:`$O`$1")
$ws.Range("H1:H1048576").Validation.Modify(3, 1, 1, "=`$AB`$1:`$AC`$1")

# Leave the selection on column L, matching where it ended up after the edit.
$null = $ws.Range("L1:L1048576").Select()
